# Fix mis-assigned higher taxon rank in the "Materials" sheet.
# For rows 2-6, the "scientificName" (AF) and "genus" (AR) columns were
# incorrectly populated with "Sagittoidea" (should be blank), and the
# "taxonRank" (AU) column was incorrectly set to "genus" (should be "class").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

for ($row = 2; $row -le 6; $row++) {
    # A bare "'" enters the cell as an (empty) text value rather than
    # clearing the cell outright, matching a blank-text genus/scientificName.
    $ws.Range("AF$row").Value = "'"
    $ws.Range("AF$row").Style = "Normal"

    $ws.Range("AR$row").Value = "'"
    $ws.Range("AR$row").Style = "Normal"

    $ws.Range("AU$row").Value = "class"
}
